$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets
$ws1.Name = "liver.nii.gz"
$ws2.Name = "tumor.nii.gz"

# --- Sheet1 (liver.nii.gz) value updates ---
$ws1.Range("B2").Value = 0.32956999999999997
$ws1.Range("C2").Value = 0.80632999999999999
$ws1.Range("D2").Value = 0.86638999999999999

$ws1.Range("B3").Value = 0.02628
$ws1.Range("C3").Value = 0.56674999999999998
$ws1.Range("D3").Value = 0.69406000000000001

$ws1.Range("B4").Value = 0.20105999999999999
$ws1.Range("C4").Value = 0.42122999999999999
$ws1.Range("D4").Value = 0.77717999999999998

$ws1.Range("B5").Value = 0.32612999999999998
$ws1.Range("C5").Value = 0.90103999999999995
$ws1.Range("D5").Value = 0.92710000000000004

$ws1.Range("C6").Value = 0.76502000000000003
$ws1.Range("D6").Value = 0.88417999999999997

$ws1.Range("B7").Value = 0.0178
$ws1.Range("C7").Value = 0.82889000000000002
$ws1.Range("D7").Value = 0.84999000000000002

$ws1.Range("B8").Value = 0.33789999999999998
$ws1.Range("C8").Value = 0.78425999999999996
$ws1.Range("D8").Value = 0.90064999999999995

$ws1.Range("B9").Value = 0.16142999999999999
$ws1.Range("C9").Value = 0.53746000000000005
$ws1.Range("D9").Value = 0.83374999999999999

$ws1.Range("B10").Value = 0.38966000000000001
$ws1.Range("C10").Value = 0.74434
$ws1.Range("D10").Value = 0.81805000000000005

$ws1.Range("B11").Value = 0.0982
$ws1.Range("C11").Value = 0.72874000000000005
$ws1.Range("D11").Value = 0.88966999999999996

$ws1.Range("C12").Value = 0.42122999999999999
$ws1.Range("D12").Value = 0.69406000000000001

$ws1.Range("B13").Value = 0.38966000000000001
$ws1.Range("C13").Value = 0.90103999999999995
$ws1.Range("D13").Value = 0.92710000000000004

$ws1.Range("B14").Value = 0.18980749999999999
$ws1.Range("C14").Value = 0.70052750000000008
$ws1.Range("D14").Value = 0.83851500000000001

$ws1.Range("B15").Value = 0.18980749999999999
$ws1.Range("C15").Value = 0.74434
$ws1.Range("D15").Value = 0.84999000000000002

# --- Sheet2 (tumor.nii.gz) value updates ---
$ws2.Range("B2").Value = 0.13816000000000001
$ws2.Range("C2").Value = 0.68960999999999995
$ws2.Range("D2").Value = 0.74992999999999999

$ws2.Range("C3").Value = 0.46988999999999997
$ws2.Range("D3").Value = 0.60265999999999997

$ws2.Range("C4").Value = 0.11624
$ws2.Range("D4").Value = 0.43231000000000003

$ws2.Range("B5").Value = 0.24221999999999999
$ws2.Range("C5").Value = 0.89371
$ws2.Range("D5").Value = 0.90003

$ws2.Range("C6").Value = 0.81328999999999996
$ws2.Range("D6").Value = 0.80223

$ws2.Range("C7").Value = 0.11824999999999999
$ws2.Range("D7").Value = 0.03136

$ws2.Range("C8").Value = 0.48979
$ws2.Range("D8").Value = 0.73196000000000006

$ws2.Range("B9").Value = 0.10535
$ws2.Range("C9").Value = 0.59228999999999998
$ws2.Range("D9").Value = 0.72536

$ws2.Range("B10").Value = 0.25814999999999999
$ws2.Range("C10").Value = 0.62980999999999998
$ws2.Range("D10").Value = 0.72579000000000005

$ws2.Range("B11").Value = 0
$ws2.Range("C11").Value = 0.44905
$ws2.Range("D11").Value = 0.63124000000000002

$ws2.Range("C12").Value = 0.11624
$ws2.Range("D12").Value = 0.03136

$ws2.Range("B13").Value = 0.25814999999999999
$ws2.Range("C13").Value = 0.89371
$ws2.Range("D13").Value = 0.90003

$ws2.Range("B14").Value = 0.083502499999999993
$ws2.Range("C14").Value = 0.52265666666666666
$ws2.Range("D14").Value = 0.60535499999999998

$ws2.Range("B15").Value = 0
$ws2.Range("C15").Value = 0.52265666666666666
$ws2.Range("D15").Value = 0.72536

# --- Column width updates (B column narrows on both sheets) ---
$ws1.Columns.Item(2).ColumnWidth = 9.5
$ws2.Columns.Item(2).ColumnWidth = 9.5

# --- Sheet view / selection updates ---
$ws1.Range("E15").Select()

# Active sheet: tumor.nii.gz becomes the tab-selected sheet
$ws2.Activate()
